# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The "sector-vab-descripcion" field (column F) is recategorized from an
# iaest-dimension to an iaest-measure. Update the metadata rows that
# describe that field accordingly and drop the mapping-file row entry,
# which only applies to dimensions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: field "kind" annotation -> now a measure, not a dimension
$ws.Range("F2").Value = "iaest-measure:sector-vab-descripcion"

# Row 3: "dim" / "medida" marker -> measure ("medida")
$ws.Range("F3").Value = "medida"

# Row 4: datatype -> measures use xsd:int (like column D "vab")
$ws.Range("F4").Value = "xsd:int"

# Row 5: mapping file reference is only used by dimensions; remove it
$ws.Range("F5").Clear()
